$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: direct assignment keeps them as text.
$ws.Range("D2").Value = '42.472.95'
$ws.Range("E2").Value = '  +1.71%  '
$ws.Range("D3").Value = '2.283.89'
$ws.Range("E3").Value = '  +0.61%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("E5").Value = '  +1.41%  '
$ws.Range("E6").Value = '  +6.17%  '
$ws.Range("E7").Value = '  +0.29%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  +2.36%  '
$ws.Range("E10").Value = '  +11.27%  '
$ws.Range("E11").Value = '  +0.35%  '
$ws.Range("E12").Value = '  -0.78%  '
$ws.Range("E13").Value = '  +0.79%  '
$ws.Range("D14").Value = '2.632.43'
$ws.Range("E14").Value = '  +0.52%  '
$ws.Range("E15").Value = '  +1.64%  '
$ws.Range("D16").Value = '2.285.77'
$ws.Range("E16").Value = '  +0.84%  '
$ws.Range("E17").Value = '  +3.44%  '
$ws.Range("D18").Value = '42.349.98'
$ws.Range("E18").Value = '  +1.69%  '
$ws.Range("E19").Value = '  +0.84%  '
$ws.Range("D20").Value = '0.0₃0913'
$ws.Range("E20").Value = '  +0.96%  '
$ws.Range("E21").Value = '  +0.96%  '
$ws.Range("E22").Value = '  +1.17%  '
$ws.Range("E23").Value = '  +0.71%  '
$ws.Range("E24").Value = '  +0.99%  '
$ws.Range("E25").Value = '  +1.69%  '
$ws.Range("E26").Value = '  -0.09%  '
$ws.Range("E27").Value = '  -0.17%  '
$ws.Range("E28").Value = '  +6.58%  '
$ws.Range("E29").Value = '  +0.28%  '
$ws.Range("E30").Value = '  +1.68%  '
$ws.Range("E31").Value = '  -0.48%  '
$ws.Range("E32").Value = '  +0.52%  '
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("E34").Value = '  +4.66%  '
$ws.Range("E35").Value = '  -0.08%  '
$ws.Range("E36").Value = '  +1.45%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("E37").Value = '  +0.65%  '
$ws.Range("B38").Value = 'WEMIXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("E38").Value = '  +0.02%  '
$ws.Range("E39").Value = '  +3.07%  '
$ws.Range("E41").Value = '  +5.68%  '
$ws.Range("E42").Value = '  +14.81%  '
$ws.Range("D43").Value = '2.000.12'
$ws.Range("E43").Value = '  -0.23%  '
$ws.Range("E44").Value = '  +1.84%  '
$ws.Range("E45").Value = '  -1.11%  '
$ws.Range("E46").Value = '  +4.24%  '
$ws.Range("E47").Value = '  -3.06%  '
$ws.Range("E48").Value = '  +1.27%  '
$ws.Range("E49").Value = '  +1.25%  '
$ws.Range("E50").Value = '  +0.41%  '
$ws.Range("E51").Value = '  +1.42%  '

# Values that parse as plain numbers (e.g. '0.999'): assigning via .Value would
# silently convert the cell to a Number. Instead, write a text-literal formula
# (="...") then Copy + PasteSpecial(values) to collapse it back to a plain text
# constant - same trick used in Excel to turn a formula into a literal value,
# without touching the cell's number format/style.
$ws.Range("D4").Formula = '="0.999"'
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4163) | Out-Null
$ws.Range("D5").Formula = '="307.76"'
$ws.Range("D5").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4163) | Out-Null
$ws.Range("D6").Formula = '="98.19"'
$ws.Range("D6").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4163) | Out-Null
$ws.Range("D7").Formula = '="0.531"'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4163) | Out-Null
$ws.Range("D9").Formula = '="0.495"'
$ws.Range("D9").Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4163) | Out-Null
$ws.Range("D10").Formula = '="36.15"'
$ws.Range("D10").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4163) | Out-Null
$ws.Range("D11").Formula = '="0.0799"'
$ws.Range("D11").Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4163) | Out-Null
$ws.Range("D13").Formula = '="6.72"'
$ws.Range("D13").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4163) | Out-Null
$ws.Range("D15").Formula = '="14.48"'
$ws.Range("D15").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4163) | Out-Null
$ws.Range("D17").Formula = '="0.800"'
$ws.Range("D17").Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4163) | Out-Null
$ws.Range("D19").Formula = '="12.62"'
$ws.Range("D19").Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4163) | Out-Null
$ws.Range("D21").Formula = '="5.99"'
$ws.Range("D21").Copy() | Out-Null
$ws.Range("D21").PasteSpecial(-4163) | Out-Null
$ws.Range("D22").Formula = '="67.80"'
$ws.Range("D22").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4163) | Out-Null
$ws.Range("D23").Formula = '="241.57"'
$ws.Range("D23").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4163) | Out-Null
$ws.Range("D25").Formula = '="1.96"'
$ws.Range("D25").Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4163) | Out-Null
$ws.Range("D27").Formula = '="23.93"'
$ws.Range("D27").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4163) | Out-Null
$ws.Range("D28").Formula = '="37.80"'
$ws.Range("D28").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4163) | Out-Null
$ws.Range("D29").Formula = '="9.56"'
$ws.Range("D29").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4163) | Out-Null
$ws.Range("D30").Formula = '="2.11"'
$ws.Range("D30").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4163) | Out-Null
$ws.Range("D31").Formula = '="159.69"'
$ws.Range("D31").Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4163) | Out-Null
$ws.Range("D32").Formula = '="5.28"'
$ws.Range("D32").Copy() | Out-Null
$ws.Range("D32").PasteSpecial(-4163) | Out-Null
$ws.Range("D34").Formula = '="3.15"'
$ws.Range("D34").Copy() | Out-Null
$ws.Range("D34").PasteSpecial(-4163) | Out-Null
$ws.Range("D35").Formula = '="0.0742"'
$ws.Range("D35").Copy() | Out-Null
$ws.Range("D35").PasteSpecial(-4163) | Out-Null
$ws.Range("D36").Formula = '="17.14"'
$ws.Range("D36").Copy() | Out-Null
$ws.Range("D36").PasteSpecial(-4163) | Out-Null
$ws.Range("D37").Formula = '="0.106"'
$ws.Range("D37").Copy() | Out-Null
$ws.Range("D37").PasteSpecial(-4163) | Out-Null
$ws.Range("D38").Formula = '="2.37"'
$ws.Range("D38").Copy() | Out-Null
$ws.Range("D38").PasteSpecial(-4163) | Out-Null
$ws.Range("D41").Formula = '="4.13"'
$ws.Range("D41").Copy() | Out-Null
$ws.Range("D41").PasteSpecial(-4163) | Out-Null
$ws.Range("D42").Formula = '="2.43"'
$ws.Range("D42").Copy() | Out-Null
$ws.Range("D42").PasteSpecial(-4163) | Out-Null
$ws.Range("D45").Formula = '="18.92"'
$ws.Range("D45").Copy() | Out-Null
$ws.Range("D45").PasteSpecial(-4163) | Out-Null
$ws.Range("D46").Formula = '="3.00"'
$ws.Range("D46").Copy() | Out-Null
$ws.Range("D46").PasteSpecial(-4163) | Out-Null
$ws.Range("D48").Formula = '="53.10"'
$ws.Range("D48").Copy() | Out-Null
$ws.Range("D48").PasteSpecial(-4163) | Out-Null
$ws.Range("D49").Formula = '="1.53"'
$ws.Range("D49").Copy() | Out-Null
$ws.Range("D49").PasteSpecial(-4163) | Out-Null
$ws.Range("D50").Formula = '="72.23"'
$ws.Range("D50").Copy() | Out-Null
$ws.Range("D50").PasteSpecial(-4163) | Out-Null
$ws.Range("D51").Formula = '="92.35"'
$ws.Range("D51").Copy() | Out-Null
$ws.Range("D51").PasteSpecial(-4163) | Out-Null
